$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09657300000000001
$ws.Range("H2").Value = 0.289719
$ws.Range("I2").Value = 0.1488791629603479
$ws.Range("J2").Value = 0.2078467214814188
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.495057333333333
$ws.Range("N2").Value = 16.485172
$ws.Range("O2").Value = 0.8161989011161211
$ws.Range("P2").Value = 0.8403205285996808
$ws.Range("Q2").Value = 0.5306741718519999
$ws.Range("R2").Value = 4.776067546668
$ws.Range("S2").Value = 0.1215150092073239
$ws.Range("T2").Value = 0.1746578668629764
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09657300000000001
$ws.Range("H3").Value = 0.289719
$ws.Range("I3").Value = 0.1488791629603479
$ws.Range("J3").Value = 0.2078467214814188
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.657666
$ws.Range("N3").Value = 1.972998
$ws.Range("O3").Value = 0.09768528951377062
$ws.Range("P3").Value = 0.1005722428790014
$ws.Range("Q3").Value = 0.063512778618
$ws.Range("R3").Value = 0.571615007562
$ws.Range("S3").Value = 0.01454330413634942
$ws.Range("T3").Value = 0.0209036109544334
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09657300000000001
$ws.Range("H4").Value = 0.289719
$ws.Range("I4").Value = 0.1488791629603479
$ws.Range("J4").Value = 0.2078467214814188
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5797745
$ws.Range("N4").Value = 1.159549
$ws.Range("O4").Value = 0.08611580937010824
$ws.Range("P4").Value = 0.0591072285213179
$ws.Range("Q4").Value = 0.0559905627885
$ws.Range("R4").Value = 0.335943376731
$ws.Range("S4").Value = 0.0128208496166746
$ws.Range("T4").Value = 0.01228524366400893
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5520940000000001
$ws.Range("H5").Value = 1.104188
$ws.Range("I5").Value = 0.8511208370396521
$ws.Range("J5").Value = 0.7921532785185812
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.495057333333333
$ws.Range("N5").Value = 16.485172
$ws.Range("O5").Value = 0.8161989011161211
$ws.Range("P5").Value = 0.8403205285996808
$ws.Range("Q5").Value = 3.033788183389333
$ws.Range("R5").Value = 18.202729100336
$ws.Range("S5").Value = 0.6946838919087973
$ws.Range("T5").Value = 0.6656626617367043
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5520940000000001
$ws.Range("H6").Value = 1.104188
$ws.Range("I6").Value = 0.8511208370396521
$ws.Range("J6").Value = 0.7921532785185812
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.657666
$ws.Range("N6").Value = 1.972998
$ws.Range("O6").Value = 0.09768528951377062
$ws.Range("P6").Value = 0.1005722428790014
$ws.Range("Q6").Value = 0.3630934526040001
$ws.Range("R6").Value = 2.178560715624001
$ws.Range("S6").Value = 0.08314198537742121
$ws.Range("T6").Value = 0.07966863192456798
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5520940000000001
$ws.Range("H7").Value = 1.104188
$ws.Range("I7").Value = 0.8511208370396521
$ws.Range("J7").Value = 0.7921532785185812
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5797745
$ws.Range("N7").Value = 1.159549
$ws.Range("O7").Value = 0.08611580937010824
$ws.Range("P7").Value = 0.0591072285213179
$ws.Range("Q7").Value = 0.320090022803
$ws.Range("R7").Value = 1.280360091212
$ws.Range("S7").Value = 0.07329495975343364
$ws.Range("T7").Value = 0.04682198485730896
